# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.25 = 34229.63 pesos`n✅ 34229.63 pesos = 8.25 = 968.09 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update N10, O10, N12, O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 121.24
$ws2.Range("O10").Value = 4150
$ws2.Range("N12").Value = 4151
$ws2.Range("O12").Value = 117.4
